$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like strings in column A to be stored as text, not auto-converted to dates.
$ws.Range("A14:A16").NumberFormat = "@"
$ws.Range("A14").Value = "2012.4.25"
$ws.Range("A15").Value = "2012.4.26"
$ws.Range("A16").Value = "2012.4.27"
$ws.Range("A14:A16").ClearFormats()

# Row 14
$ws.Range("B14").Value = "将新素材加入纹理图"
$ws.Range("D14").Value = 2.5
$ws.Range("B14").WrapText = $true

# Row 15
$ws.Range("B15").Value = "1.加入飞行动画。2.调试场景切换bug，未成功。3.修改一处操作，要先触摸到飞行动物才能移动。"
$ws.Range("C15").Value = "场景切换，body释放会段错，现在暂时没有释放。"
$ws.Range("D15").Value = 3
$ws.Range("B15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 40.5

# Row 16
$ws.Range("B16").Value = "初步完成飞行动物的动画"
$ws.Range("C16").Value = "精灵在entity里申明的，显示不出来，必须要在子类再次申明一下。没明白为什么。`n飞行动物碰撞边框后，方向有点错误"
$ws.Range("D16").Value = 3
$ws.Range("B16").WrapText = $true
$ws.Range("C16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 40.5

# Update sheet view to match the target (topLeftCell A7, selection B18)
$ws.Application.ActiveWindow.ScrollRow = 7
$null = $ws.Range("B18").Select()
